# queries - ranking - rec
# Re-rank tied UFs ("qtd" sheet) and refresh the success-rate ("tx-sucesso")
# ranking so the row order matches the recalculated source queries.

$wb = $excel.ActiveWorkbook

# --- Sheet "qtd": re-order UF labels among tied counts -------------------
$wsQtd = $wb.Worksheets.Item("qtd")
$wsQtd.Range("A10").Value = "PB"
$wsQtd.Range("A11").Value = "BA"
$wsQtd.Range("A17").Value = "ES"
$wsQtd.Range("A18").Value = "MS"
$wsQtd.Range("A20").Value = "RO"
$wsQtd.Range("A21").Value = "AM"
$wsQtd.Range("A22").Value = "MA"
$wsQtd.Range("A23").Value = "SE"
$wsQtd.Range("A24").Value = "MT"

# --- Sheet "tx-sucesso": re-order the 100%-tie group ----------------------
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A2").Value = "RO"
$wsTx.Range("A3").Value = ""
$wsTx.Range("A5").Value = "MA"
